# Insert a new weekly price record for "Vega Modelo de Temuco - Mango"
# right before the current row 567, shifting all subsequent rows down by
# one (old row 567 -> 568, ..., old row 615 -> 616).
#
# The new row duplicates every column from the row that ends up directly
# below it (the former row 567) except for the fields that actually carry
# new data for this record: Fecha (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 567:615 down to 568:616, leaving a blank row 567 behind.
$ws.Rows.Item(567).Insert()

# Seed the new row with the same data as the row that is now right below
# it (the former row 567), so every column besides the ones we overwrite
# next keeps matching formatting/values.
$ws.Range("A568:T568").Copy($ws.Range("A567:T567"))

# Overwrite with this week's actual figures.
$ws.Range("D567").Value = 45132
$ws.Range("M567").Value = 1000
$ws.Range("N567").Value = 8000
$ws.Range("O567").Value = 8000
$ws.Range("P567").Value = 8000
$ws.Range("S567").Value = 2000
